# Update background dataset labels (casing fix) and refresh the active
# selection, per commit "fix: update background datatsets and makefile".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 15 / Row 16 held the two material names that need their casing
# normalized ("Title Case" -> lower/sentence case), keeping all other
# data (wastage % in column B, notes in column C) untouched.
$ws.Range("A15").Value = '6" galvanized steel stud framing'
$ws.Range("A16").Value = "Gypsum wall board"

# Update the saved selection/active cell for the sheet view.
$ws.Range("E13").Select()
